$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Mean (H) and Std (I) values for PreferenceOrder prediction-type rows
# across all Base_Learner blocks (RF, XGBoost, ETC, LightGBM).
$ws.Range("H26").Value = 0.63183
$ws.Range("I26").Value = 0.03408
$ws.Range("H27").Value = 0.06578000000000001
$ws.Range("I27").Value = 0.03469
$ws.Range("H28").Value = 0.63624
$ws.Range("I28").Value = 0.03366
$ws.Range("H29").Value = 0.09387
$ws.Range("I29").Value = 0.04478
$ws.Range("H30").Value = 0.63331
$ws.Range("I30").Value = 0.03423
$ws.Range("H31").Value = 0.06481000000000001
$ws.Range("I31").Value = 0.03269
$ws.Range("H32").Value = 0.63844
$ws.Range("I32").Value = 0.03159
$ws.Range("H33").Value = 0.09872
$ws.Range("I33").Value = 0.04471
$ws.Range("H34").Value = 0.63809
$ws.Range("I34").Value = 0.02424
$ws.Range("H35").Value = 0.01259
$ws.Range("I35").Value = 0.0196
$ws.Range("H36").Value = 0.63803
$ws.Range("I36").Value = 0.02447
$ws.Range("H37").Value = 0.01259
$ws.Range("I37").Value = 0.0196
$ws.Range("H38").Value = 0.63835
$ws.Range("I38").Value = 0.02389
$ws.Range("H39").Value = 0.01259
$ws.Range("I39").Value = 0.0196
$ws.Range("H40").Value = 0.6379
$ws.Range("I40").Value = 0.02414
$ws.Range("H41").Value = 0.01259
$ws.Range("I41").Value = 0.0196
$ws.Range("H66").Value = 0.56007
$ws.Range("I66").Value = 0.0274
$ws.Range("H67").Value = 0.0164
$ws.Range("I67").Value = 0.01898
$ws.Range("H68").Value = 0.57987
$ws.Range("I68").Value = 0.0281
$ws.Range("H69").Value = 0.04058
$ws.Range("I69").Value = 0.03023
$ws.Range("H70").Value = 0.56231
$ws.Range("I70").Value = 0.02897
$ws.Range("H71").Value = 0.01545
$ws.Range("I71").Value = 0.01918
$ws.Range("H72").Value = 0.58256
$ws.Range("I72").Value = 0.02948
$ws.Range("H73").Value = 0.0502
$ws.Range("I73").Value = 0.03026
$ws.Range("H74").Value = 0.58085
$ws.Range("I74").Value = 0.02
$ws.Range("H75").Value = 0.00771
$ws.Range("I75").Value = 0.01481
$ws.Range("H76").Value = 0.58777
$ws.Range("I76").Value = 0.01904
$ws.Range("H77").Value = 0.01159
$ws.Range("I77").Value = 0.01692
$ws.Range("H78").Value = 0.5805900000000001
$ws.Range("I78").Value = 0.0202
$ws.Range("H79").Value = 0.00771
$ws.Range("I79").Value = 0.01481
$ws.Range("H80").Value = 0.5879
$ws.Range("I80").Value = 0.01878
$ws.Range("H81").Value = 0.01159
$ws.Range("I81").Value = 0.01692
$ws.Range("H106").Value = 0.63146
$ws.Range("I106").Value = 0.03553
$ws.Range("H107").Value = 0.05312
$ws.Range("I107").Value = 0.04643
$ws.Range("H108").Value = 0.63363
$ws.Range("I108").Value = 0.03592
$ws.Range("H109").Value = 0.07055
$ws.Range("I109").Value = 0.06122
$ws.Range("H110").Value = 0.63204
$ws.Range("I110").Value = 0.03524
$ws.Range("H111").Value = 0.05312
$ws.Range("I111").Value = 0.04694
$ws.Range("H112").Value = 0.63537
$ws.Range("I112").Value = 0.0387
$ws.Range("H113").Value = 0.07826
$ws.Range("I113").Value = 0.06501999999999999
$ws.Range("H114").Value = 0.63611
$ws.Range("I114").Value = 0.02588
$ws.Range("H115").Value = 0.01261
$ws.Range("I115").Value = 0.02585
$ws.Range("H116").Value = 0.63605
$ws.Range("I116").Value = 0.02609
$ws.Range("H117").Value = 0.01261
$ws.Range("I117").Value = 0.02585
$ws.Range("H118").Value = 0.63611
$ws.Range("I118").Value = 0.02588
$ws.Range("H119").Value = 0.01261
$ws.Range("I119").Value = 0.02585
$ws.Range("H120").Value = 0.63605
$ws.Range("I120").Value = 0.02609
$ws.Range("H121").Value = 0.01261
$ws.Range("I121").Value = 0.02585
$ws.Range("H146").Value = 0.5944199999999999
$ws.Range("I146").Value = 0.03104
$ws.Range("H147").Value = 0.0174
$ws.Range("I147").Value = 0.01868
$ws.Range("H148").Value = 0.60934
$ws.Range("I148").Value = 0.03316
$ws.Range("H149").Value = 0.05131
$ws.Range("I149").Value = 0.04891
$ws.Range("H150").Value = 0.5948
$ws.Range("I150").Value = 0.03055
$ws.Range("H151").Value = 0.01837
$ws.Range("I151").Value = 0.01838
$ws.Range("H152").Value = 0.61322
$ws.Range("I152").Value = 0.03308
$ws.Range("H153").Value = 0.06393
$ws.Range("I153").Value = 0.04687
$ws.Range("H154").Value = 0.61072
$ws.Range("I154").Value = 0.02926
$ws.Range("H155").Value = 0.01542
$ws.Range("I155").Value = 0.02352
$ws.Range("H156").Value = 0.61465
$ws.Range("I156").Value = 0.03107
$ws.Range("H157").Value = 0.02316
$ws.Range("I157").Value = 0.02919
$ws.Range("H158").Value = 0.6107900000000001
$ws.Range("I158").Value = 0.02931
$ws.Range("H159").Value = 0.01542
$ws.Range("I159").Value = 0.02352
$ws.Range("H160").Value = 0.6149
$ws.Range("I160").Value = 0.03114
$ws.Range("H161").Value = 0.02314
$ws.Range("I161").Value = 0.0314
